$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Reorder the header labels in row 1 (A1:F1) to match the finished
# input-list generation order.
$ws.Range("A1").Value = "living_rooms_1"
$ws.Range("B1").Value = "bedrooms_1"
$ws.Range("C1").Value = "bedrooms_2"
$ws.Range("D1").Value = "kitchens_1"
$ws.Range("E1").Value = "living_rooms_2"
$ws.Range("F1").Value = "kitchens_2"
